# Fixed headers in test data
# The "Connections" sheet had header labels OriginChain / OriginProcess /
# DestinationChain that no longer matched the naming convention used
# elsewhere. Rename them to Origin_Chain / Origin_Unit / Destination_Chain.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Connections")

$ws.Range("A1").Value = "Origin_Chain"
$ws.Range("B1").Value = "Origin_Unit"
$ws.Range("F1").Value = "Destination_Chain"

# Restore the user's last selection on this sheet.
$ws.Range("C9").Select()
